# Update Name of Algo
# Applies updated imputed values produced by the (re-named) KNN algorithm
# to the result_data_KNN.xlsx workbook. Only specific numeric cells in
# columns A, C and D change value; everything else (layout, headers,
# formatting) stays the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D3").Value = -7.702
$ws.Range("D4").Value = -7.892
$ws.Range("C7").Value = -12.983
$ws.Range("A8").Value = -22.273
$ws.Range("A10").Value = -21.82
$ws.Range("D11").Value = -7.2
$ws.Range("A12").Value = -21.594
$ws.Range("C14").Value = -12.996
$ws.Range("D14").Value = -7.986
$ws.Range("C15").Value = -13.745
$ws.Range("A18").Value = -21.841
$ws.Range("C18").Value = -12.098
$ws.Range("D18").Value = -7.867999999999999
$ws.Range("D19").Value = -8.254999999999999
$ws.Range("C20").Value = -12.684
$ws.Range("D21").Value = -8.353999999999999
$ws.Range("A25").Value = -21.862
$ws.Range("D27").Value = -8.568999999999999
$ws.Range("C29").Value = -12.181
$ws.Range("C30").Value = -12.421
$ws.Range("C31").Value = -13.363
$ws.Range("D31").Value = -8.425999999999998
$ws.Range("C35").Value = -12.518
$ws.Range("A37").Value = -20.044
$ws.Range("D38").Value = -7.813
$ws.Range("C40").Value = -12.782
$ws.Range("D42").Value = -8.303999999999998
$ws.Range("C44").Value = -12.663
$ws.Range("D44").Value = -7.564
$ws.Range("D47").Value = -7.395
$ws.Range("C50").Value = -13.542
$ws.Range("C54").Value = -12.524
$ws.Range("A55").Value = -21.894
$ws.Range("D56").Value = -8.164
$ws.Range("D58").Value = -8.401
$ws.Range("D65").Value = -7.437
$ws.Range("A68").Value = -21.593
$ws.Range("D73").Value = -8.193000000000001
$ws.Range("C76").Value = -13.314
$ws.Range("A77").Value = -20.42
$ws.Range("A78").Value = -20.099
$ws.Range("A79").Value = -21.771
$ws.Range("A80").Value = -20.228
$ws.Range("A81").Value = -21.759
$ws.Range("A82").Value = -22.124
$ws.Range("A84").Value = -22.177
$ws.Range("C87").Value = -13.215
$ws.Range("C88").Value = -13.085
$ws.Range("D90").Value = -7.452
$ws.Range("C92").Value = -11.474
$ws.Range("D92").Value = -6.709000000000001
$ws.Range("D94").Value = -7.129
$ws.Range("D95").Value = -7.836000000000001
$ws.Range("C96").Value = -12.665
$ws.Range("C98").Value = -13.446
$ws.Range("A101").Value = -20.846
$ws.Range("C101").Value = -12.585
$ws.Range("D101").Value = -7.816
$ws.Range("A102").Value = -19.812
$ws.Range("C102").Value = -12.567
